# Updates the crypto price/volume table (columns D "Price" and E
# "Volume(1h)", rows 2-51) on Sheet1 to the refreshed values captured by
# the latest GitHub Actions run, per the commit "Updated cryptos list on
# Sat Apr  1 09:55:51 UTC 2023 with GitHub Actions".
#
# All of these source cells are plain text (inline strings) — prices like
# "28.480.62" use '.' as a thousands separator (not a decimal point), so
# everything must stay text. Two of the new Price values (D40, D45) are
# otherwise indistinguishable from genuine numbers and would have a
# significant trailing zero silently stripped by Excel's normal
# type-inference (e.g. "0.6250" -> 0.625), so those two cells are pinned
# to the Text number format before the value is written, exactly as you'd
# do by hand in Excel to stop that auto-conversion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.480.62'
$ws.Range('E2').Value = '  +2.38%  '

$ws.Range('D3').Value = '1.828.17'
$ws.Range('E3').Value = '  +2.14%  '

$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').Value = '315.69'
$ws.Range('E5').Value = '  -0.24%  '

$ws.Range('E6').Value = '  +0.04%  '

$ws.Range('D7').Value = '0.5072'
$ws.Range('E7').Value = '  -5.08%  '

$ws.Range('D8').Value = '0.3913'
$ws.Range('E8').Value = '  +1.75%  '

$ws.Range('D9').Value = '0.07698'
$ws.Range('E9').Value = '  +3.79%  '

$ws.Range('D10').Value = '41.94'
$ws.Range('E10').Value = '  +1.56%  '

$ws.Range('D11').Value = '1.114'
$ws.Range('E11').Value = '  +2.94%  '

$ws.Range('D12').Value = '21.06'
$ws.Range('E12').Value = '  +3.97%  '

$ws.Range('D13').Value = '6.274'
$ws.Range('E13').Value = '  +1.71%  '

$ws.Range('D14').Value = '1.002'
$ws.Range('E14').Value = '  +0.15%  '

$ws.Range('D15').Value = '7.562'
$ws.Range('E15').Value = '  +1.99%  '

$ws.Range('D16').Value = '1.824.87'
$ws.Range('E16').Value = '  +1.98%  '

$ws.Range('D17').Value = '93.48'
$ws.Range('E17').Value = '  +6.19%  '

$ws.Range('D18').Value = '0.00001086'
$ws.Range('E18').Value = '  +2.80%  '

$ws.Range('D19').Value = '0.06662'
$ws.Range('E19').Value = '  +2.18%  '

$ws.Range('E20').Value = '  +3.13%  '

$ws.Range('E21').Value = '  +0.10%  '

$ws.Range('D22').Value = '6.166'
$ws.Range('E22').Value = '  +3.77%  '

$ws.Range('D23').Value = '28.511.90'
$ws.Range('E23').Value = '  +2.29%  '

$ws.Range('E24').Value = '  +0.05%  '

$ws.Range('D25').Value = '2.258'
$ws.Range('E25').Value = '  +7.87%  '

$ws.Range('D26').Value = '156.99'
$ws.Range('E26').Value = '  +0.11%  '

$ws.Range('D27').Value = '20.63'
$ws.Range('E27').Value = '  +2.69%  '

$ws.Range('D28').Value = '2.038.15'
$ws.Range('E28').Value = '  +2.11%  '

$ws.Range('E29').Value = '  +5.55%  '

$ws.Range('D30').Value = '125.31'
$ws.Range('E30').Value = '  +3.37%  '

$ws.Range('D31').Value = '1.133'
$ws.Range('E31').Value = '  +3.65%  '

$ws.Range('D32').Value = '0.1089'
$ws.Range('E32').Value = '  -0.41%  '

$ws.Range('D33').Value = '5.675'
$ws.Range('E33').Value = '  +3.55%  '

$ws.Range('D34').Value = '3.662'
$ws.Range('E34').Value = '  +0.41%  '

$ws.Range('D35').Value = '0.07086'
$ws.Range('E35').Value = '  +1.90%  '

$ws.Range('D36').Value = '0.2225'
$ws.Range('E36').Value = '  +1.44%  '

$ws.Range('D37').Value = '0.02327'
$ws.Range('E37').Value = '  +3.01%  '

$ws.Range('D38').Value = '8.887'
$ws.Range('E38').Value = '  +6.46%  '

$ws.Range('D39').Value = '5.156'
$ws.Range('E39').Value = '  +2.19%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6250'
$ws.Range('E40').Value = '  +2.84%  '

$ws.Range('E41').Value = '  -0.18%  '

$ws.Range('D42').Value = '1.188'
$ws.Range('E42').Value = '  +1.30%  '

$ws.Range('E43').Value = '  +0.05%  '

$ws.Range('D44').Value = '1.398'
$ws.Range('E44').Value = '  -1.18%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.50'
$ws.Range('E45').Value = '  +2.03%  '

$ws.Range('D46').Value = '0.5914'
$ws.Range('E46').Value = '  +4.17%  '

$ws.Range('D47').Value = '3.718'
$ws.Range('E47').Value = '  +1.04%  '

$ws.Range('D48').Value = '124.94'
$ws.Range('E48').Value = '  +0.55%  '

$ws.Range('E49').Value = '  +3.96%  '

$ws.Range('D50').Value = '1.194'
$ws.Range('E50').Value = '  +2.29%  '

$ws.Range('D51').Value = '0.06919'
$ws.Range('E51').Value = '  +1.92%  '
